# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2310
#   *_new -> *_FV2404
# Then freeze the header row and wrap the data range in a native Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header-row cells (row 1) so formatversion suffixes replace
#    the old "_old"/"_new" naming scheme.
# ---------------------------------------------------------------------------
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -like "*_old") {
        $cell.Value = ($v -replace "_old$", "_FV2310")
    } elseif ($v -like "*_new") {
        $cell.Value = ($v -replace "_new$", "_FV2404")
    }
}

# ---------------------------------------------------------------------------
# 2. Freeze the top (header) row.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3. Turn the used range into a formatted Table (ListObject) with the
#    (now renamed) header row, matching the source data range A1:U56.
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:U56")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

Write-Output "Header renaming, frozen header row and Table1 applied."
